# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# worksheet with refreshed figures, as produced by the scheduled
# GitHub Actions crypto-list refresh job.
#
# Note: several Price values look like plain decimal numbers (e.g.
# "0.9390", "20.87"); Excel would otherwise auto-convert these to
# numeric values (losing significant trailing zeros / the original
# text type). To preserve them as text exactly as authored, the
# Range is pre-formatted as Text ("@") before assigning such values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "28.200.95"
$ws.Range("E2").Value = "  +2.46%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.916.83"
$ws.Range("E3").Value = "  +2.11%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.90%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.52"
$ws.Range("E5").Value = "  +1.29%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  -0.79%  "

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4843"
$ws.Range("E7").Value = "  +0.98%  "

# Row 8 - Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3831"
$ws.Range("E8").Value = "  +1.40%  "

# Row 9 - Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07385"
$ws.Range("E9").Value = "  -0.06%  "

# Row 10 - Polygon
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9390"
$ws.Range("E10").Value = "  -0.22%  "

# Row 11 - Solana
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.87"
$ws.Range("E11").Value = "  +0.54%  "

# Row 12 - TRON
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07808"
$ws.Range("E12").Value = "  -0.79%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.926.00"
$ws.Range("E13").Value = "  +2.52%  "

# Row 14 - Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.511"
$ws.Range("E14").Value = "  +1.10%  "

# Row 15 - Chainlink
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.648"
$ws.Range("E15").Value = "  +0.55%  "

# Row 16 - Litecoin
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.26"
$ws.Range("E16").Value = "  +0.07%  "

# Row 17 - BinanceUSD
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.008"
$ws.Range("E17").Value = "  -0.76%  "

# Row 18 - ShibaInu
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008847"
$ws.Range("E18").Value = "  -1.64%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  -0.80%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "28.223.90"
$ws.Range("E20").Value = "  +2.48%  "

# Row 21 - Avalanche
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.87"
$ws.Range("E21").Value = "  -0.66%  "

# Row 22 - Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.165"
$ws.Range("E22").Value = "  +0.31%  "

# Row 23 - WrappedliquidstakedEther2.0
$ws.Range("D23").Value = "2.183.51"
$ws.Range("E23").Value = "  +3.62%  "

# Row 24 - Cosmos
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.93"
$ws.Range("E24").Value = "  +1.95%  "

# Row 25 - Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.65"
$ws.Range("E25").Value = "  +1.75%  "

# Row 26 - Toncoin
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.921"
$ws.Range("E26").Value = "  -1.51%  "

# Row 27 - EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.57"
$ws.Range("E27").Value = "  -0.10%  "

# Row 28 - LidoDAOToken
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.106"
$ws.Range("E28").Value = "  +3.83%  "

# Row 29 - BitcoinCash
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "116.35"
$ws.Range("E29").Value = "  +0.11%  "

# Row 30 - InternetComputer(DFINITY)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.971"
$ws.Range("E30").Value = "  -0.67%  "

# Row 31 - Stellar
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08915"
$ws.Range("E31").Value = "  -0.12%  "

# Row 32 - HuobiToken
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.362"
$ws.Range("E32").Value = "  +1.11%  "

# Row 33 - ARBITRUM
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.254"
$ws.Range("E33").Value = "  +2.76%  "

# Row 34 - ImmutableX
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7733"
$ws.Range("E34").Value = "  +2.70%  "

# Row 35 - Filecoin
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.699"
$ws.Range("E35").Value = "  +1.91%  "

# Row 36 - RenderToken
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.657"
$ws.Range("E36").Value = "  -1.83%  "

# Row 37 - VeChain
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02054"
$ws.Range("E37").Value = "  -1.17%  "

# Row 38 - TrustWalletToken
$ws.Range("E38").Value = "  -1.81%  "

# Row 39 - Hedera
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05332"
$ws.Range("E39").Value = "  +0.35%  "

# Row 40 - TheSandbox
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5529"
$ws.Range("E40").Value = "  +2.81%  "

# Row 41 - MXToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.004"
$ws.Range("E41").Value = "  -0.13%  "

# Row 42 - FraxShare
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.074"
$ws.Range("E42").Value = "  -0.43%  "

# Row 43 - Algorand
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1532"
$ws.Range("E43").Value = "  +0.30%  "

# Row 44 - Aptos
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.469"
$ws.Range("E44").Value = "  +0.15%  "

# Row 45 - EnergySwap
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.75"
$ws.Range("E45").Value = "  +0.37%  "

# Row 46 - Decentraland
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4863"
$ws.Range("E46").Value = "  +0.28%  "

# Row 47 - Quant
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "107.23"
$ws.Range("E47").Value = "  +3.96%  "

# Row 49 - NEARProtocol
$ws.Range("E49").Value = "  -0.37%  "

# Row 50 - Aave
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.81"
$ws.Range("E50").Value = "  +2.30%  "

# Row 51 - Cronos
$ws.Range("E51").Value = "  +0.01%  "
